$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shorten the room labels in column D from "Room NNN" to "R NNN"
$rng = $ws.Range("D1:D84")
$rng.Replace("Room ", "R ") | Out-Null

# Leave the selection on column E, matching the author's final UI state
$ws.Range("E1:E1048576").Select() | Out-Null
